$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename column header B1 from "errors" to "amount"
$ws.Range("B1").Value = "amount"

# 2. Append a new data row (row 4) with a plain (non-error) pair of numbers,
#    formatted like the existing data cells (Calibri font, centered) but
#    without the outer-table border used by the header/data block above it.
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 5

# Borrow the font/number-format/alignment from an existing data cell so the
# new row visually matches the rest of the table...
$ws.Range("B3").Copy()
$ws.Range("A4:B4").PasteSpecial(-4122)

# ...then strip the border that the copied cell had, since the new row sits
# outside the bordered block.
$ws.Range("A4:B4").Borders.LineStyle = -4142
